$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "УТЗ, ЦВД"
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 1

$ws.Range("D9").Value = 8.767939928491167
$ws.Range("E9").Value = 0.02834233223622572
$ws.Range("F9").Value = 0.0009996105588871713

$ws.Range("C14").Value = 1.998035383653871
$ws.Range("D14").Value = 9.767300492850904
$ws.Range("F14").Value = 0.0009991108376305076

$ws.Range("C18").Value = -0.999360564359737

$ws.Range("C22").Value = -14126.68016159745

$ws.Range("H24").Value = 28833.51946702796

$ws.Range("F27").Value = 0.95

$ws.Range("C32").Value = 2.011978024246901
$ws.Range("D32").Value = 8.767939928491167
$ws.Range("E32").Value = 0.02852407789332977
$ws.Range("F32").Value = 0.0009996352088091508
